# Update the "Förändrad" (Changed) date column (C) for all data rows
# from the old value (45190 -> 2023-09-21) to the new value
# (45192 -> 2023-09-23), matching the committed workbook change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 339

$ws.Range("C$firstRow`:C$lastRow").Value = 45192
